$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: Level Walk + Slope up and down (LC-style wrapped accuracy block) ---
$ws.Cells.Item(18, 1).Value = "Level Walk + Slope up and down: 495 subjects"
$ws.Cells.Item(18, 2).Value = "prediction accuracy on test set: 65.3137%`nprediction accuracy on test set: 71.5867%`nprediction accuracy on test set: 69.7417%`nprediction accuracy on test set: 69.0037%`nprediction accuracy on test set: 70.8487%`nprediction accuracy on test set: 66.0517%`nprediction accuracy on test set: 69.7417%`nprediction accuracy on test set: 64.2066%`nprediction accuracy on test set: 68.2657%`nprediction accuracy on test set: 68.2657%`n"
$ws.Cells.Item(18, 2).WrapText = $true
$ws.Cells.Item(18, 3).Value = 677159
$ws.Cells.Item(18, 4).Value = 1354
$ws.Cells.Item(18, 5).Value = "Wavelet + RF"
$ws.Rows.Item(18).RowHeight = 78

# --- Row 19: val_accuracy summary row ---
$ws.Cells.Item(19, 1).Value = "Level Walk + Slope up and down: 495 subjects"
$ws.Cells.Item(19, 2).Value = "val_accuracy: 0.7085"
$ws.Cells.Item(19, 3).Value = 677159
$ws.Cells.Item(19, 4).Value = 1354
$ws.Cells.Item(19, 5).Value = "Wavelet + CNN"
$ws.Rows.Item(19).RowHeight = 40.8

# --- Column A width widened to fit the new long label ---
$ws.Columns.Item(1).ColumnWidth = 44.3

# --- View state: selection moves to the newly added B19 cell ---
$ws.Range("B19").Select() | Out-Null
